# Applies the cryptos list price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text value.
# Values are written with a leading apostrophe and Value2 so that numeric-looking
# strings (e.g. "214.54", "29.933.59") are stored as text, matching the source data,
# and the cell style is reset afterwards so no stray "quote prefix" style is introduced.
$updates = @(
    @("D2", "29.933.59"),
    @("E2", "  +0.84%  "),
    @("D3", "1.632.99"),
    @("E3", "  +1.74%  "),
    @("D4", "0.999"),
    @("E4", "  -0.19%  "),
    @("D5", "214.54"),
    @("E5", "  +0.94%  "),
    @("E6", "  +1.26%  "),
    @("E7", "  -0.14%  "),
    @("D8", "29.37"),
    @("E8", "  +9.17%  "),
    @("D9", "0.259"),
    @("E9", "  +2.93%  "),
    @("D10", "0.0614"),
    @("E10", "  +2.32%  "),
    @("D11", "0.0917"),
    @("E11", "  +1.00%  "),
    @("D12", "1.866.91"),
    @("E12", "  +1.74%  "),
    @("D13", "1.629.57"),
    @("E13", "  +1.54%  "),
    @("D14", "0.571"),
    @("E14", "  +6.24%  "),
    @("D15", "3.93"),
    @("E15", "  +5.11%  "),
    @("D16", "29.984.51"),
    @("D17", "9.02"),
    @("E17", "  +18.99%  "),
    @("D18", "65.02"),
    @("E18", "  +2.38%  "),
    @("D19", "245.10"),
    @("E19", "  +1.72%  "),
    @("E20", "  +2.17%  "),
    @("D21", "0.999"),
    @("E21", "  -0.16%  "),
    @("E22", "  +3.67%  "),
    @("D23", "9.62"),
    @("E23", "  +4.25%  "),
    @("D24", "2.13"),
    @("E24", "  +2.34%  "),
    @("D25", "158.20"),
    @("E25", "  +2.13%  "),
    @("D26", "15.75"),
    @("E26", "  +1.63%  "),
    @("E27", "  +3.20%  "),
    @("D28", "6.64"),
    @("E28", "  +3.50%  "),
    @("E29", "  -0.20%  "),
    @("E30", "  +3.20%  "),
    @("E31", "  +5.18%  "),
    @("D32", "3.37"),
    @("E32", "  +4.71%  "),
    @("E33", "  +2.69%  "),
    @("D34", "1.429.61"),
    @("E34", "  -0.08%  "),
    @("D35", "1.65"),
    @("E35", "  +6.86%  "),
    @("D36", "1.04"),
    @("E36", "  +1.32%  "),
    @("D37", "2.88"),
    @("E37", "  +2.22%  "),
    @("D38", "2.29"),
    @("E38", "  -1.24%  "),
    @("D39", "0.0171"),
    @("E39", "  +3.42%  "),
    @("D40", "0.557"),
    @("E40", "  +4.11%  "),
    @("B41", "ARBITRUM"),
    @("C41", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"),
    @("D41", "0.835"),
    @("E41", "  +4.28%  "),
    @("B42", "Kaspa"),
    @("C42", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"),
    @("D42", "0.0500"),
    @("E42", "  +1.79%  "),
    @("B43", "Aave"),
    @("C43", "https://coinranking.com/coin/ixgUfzmLR+aave-aave"),
    @("D43", "71.33"),
    @("E43", "  +8.58%  "),
    @("B44", "RenderToken"),
    @("C44", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"),
    @("D44", "1.97"),
    @("E44", "  +0.90%  "),
    @("B45", "BitcoinSV"),
    @("C45", "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"),
    @("D45", "54.70"),
    @("E45", "  +1.12%  "),
    @("B46", "WEMIXToken"),
    @("C46", "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"),
    @("D46", "1.02"),
    @("E46", "  +4.74%  "),
    @("B47", "PaxDollar"),
    @("C47", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"),
    @("D47", "0.999"),
    @("E47", "  -0.14%  "),
    @("D48", "5.43"),
    @("E48", "  +2.81%  "),
    @("D49", "1.774.43"),
    @("E49", "  +1.73%  "),
    @("D50", "89.33"),
    @("E50", "  +3.23%  "),
    @("D51", "0.0₆0108"),
    @("E51", "  +3.62%  ")
)

foreach ($update in $updates) {
    $cellRef = $update[0]
    $newValue = $update[1]
    $cell = $ws.Range($cellRef)
    $cell.Value2 = "'" + $newValue
    $cell.Style = "Normal"
}
